$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "3done"
$ws.Range("D16").Value = "3done"
$ws.Range("D24").Value = "3done"
$ws.Range("D26").Value = "3done"
$ws.Range("D29").Value = "3done"
